$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20: heater / space_heater / 1800
$ws.Range("A20").Value = "heater"
$ws.Range("B20").Value = "space_heater"
$ws.Range("C20").Value = 1800

# Row 21: iron / 1000
$ws.Range("A21").Value = "iron"
$ws.Range("C21").Value = 1000

# Row 22: blender / 350
$ws.Range("A22").Value = "blender"
$ws.Range("C22").Value = 350

# Row 23: freezer / 440
$ws.Range("A23").Value = "freezer"
$ws.Range("C23").Value = 440

# Row 24: monitor / 200
$ws.Range("A24").Value = "monitor"
$ws.Range("C24").Value = 200

# Row 25: echo / alexa / 3
$ws.Range("A25").Value = "echo"
$ws.Range("B25").Value = "alexa"
$ws.Range("C25").Value = 3

# Update the active selection to match the saved view state
$ws.Range("A26").Select()
